$d = $word.ActiveDocument

$replacements = @(
    @("2024-03-04 Monday", "2024-03-05 Tuesday"),
    @("42÷7=", "13÷5="),
    @("22÷8=", "85÷3="),
    @("37÷3=", "13÷7="),
    @("51÷9=", "23÷2="),
    @("21÷6=", "27÷9="),
    @("75÷6=", "45÷6="),
    @("90÷3=", "49÷5="),
    @("97÷7=", "26÷9="),
    @("89÷5=", "11÷8="),
    @("22÷6=", "72÷7="),
    @("10÷4=", "70÷9="),
    @("68÷7=", "79÷4="),
    @("56÷2=", "37÷8="),
    @("93÷2=", "29÷7="),
    @("84÷4=", "88÷6="),
    @("30÷4=", "69÷6="),
    @("82÷8=", "75÷5="),
    @("50÷6=", "80÷9="),
    @("92÷7=", "19÷4="),
    @("62÷2=", "40÷4="),
    @("58÷5=", "50÷8="),
    @("65÷2=", "28÷3="),
    @("41÷7=", "43÷6="),
    @("39÷7=", "70÷6="),
    @("31÷5=", "53÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
